# Script 1 - atualização automática de dados
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = "Atividades imobiliárias"; B = 28.31208085602502; C = "2014 / 2023" },
    @{ Row = 3;  A = "Atividades financeiras, de seguros e serviços relacionados"; B = 23.65397356444598; C = "2014 / 2023" },
    @{ Row = 4;  A = "Eletricidade e gás, água, esgoto, atividades de gestão de resíduos e descontaminação"; B = 22.70215620370851; C = "2014 / 2023" },
    @{ Row = 5;  A = "Informação e comunicação"; B = 18.69747267971518; C = "2014 / 2023" },
    @{ Row = 6;  A = "Agropecuária"; B = 5.682545327113388; C = "2014 / 2023" },
    @{ Row = 7;  A = "Administração, defesa, educação e saúde públicas e seguridade social"; B = 2.548568009332001; C = "2014 / 2023" },
    @{ Row = 8;  A = "Indústrias extrativas"; B = 69.71390886407463; C = "2022 / 2023" },
    @{ Row = 9;  A = "Agropecuária"; B = 7.64923808343667; C = "2022 / 2023" },
    @{ Row = 10; A = "Informação e comunicação"; B = 7.133426276710783; C = "2022 / 2023" },
    @{ Row = 11; A = "Atividades financeiras, de seguros e serviços relacionados"; B = 4.80083445372183; C = "2022 / 2023" },
    @{ Row = 12; A = "Comércio e reparação de veículos automotores e motocicletas"; B = 2.010492569309139; C = "2022 / 2023" },
    @{ Row = 13; A = "Indústrias de transformação"; B = 1.925088415038929; C = "2022 / 2023" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
}
